$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 800
$ws.Range("B3").Value = 900
$ws.Range("B4").Value = 300
$ws.Range("B5").Value = 150
$ws.Range("B6").Value = 75
$ws.Range("B7").Value = 235
$ws.Range("B8").Value = 127
